$orderData = @(
    @("100091365", "AVR2109160159882"),
    @("100091366", "AVC2109160159884"),
    @("100091368", "AVC2109160159888"),
    @("100091371", "AVR2109160159892"),
    @("100091373", "AVR2109160159895"),
    @("100091379", "AVC2109160159911"),
    @("100091380", "AVC2109160159912"),
    @("100091381", "AVC2109160159913"),
    @("100091382", "AVC2109160159914"),
    @("100091383", "AVC2109160159915"),
    @("100091385", "AVC2109160159922"),
    @("100091386", "AVC2109160159923"),
    @("100091387", "AVC2109160159924"),
    @("100091388", "AVR2109160159925"),
    @("100091390", "AVR2109160159930"),
    @("100091391", "AVR2109160159933"),
    @("100091392", "AVR2109160159935"),
    @("100091393", "AVR2109160159937"),
    @("100091394", "AVR2109160159939"),
    @("100091395", "AVR2109160159941"),
    @("100091396", "AVR2109160159943"),
    @("100091397", "AVR2109160159945"),
    @("100091399", "AVR2109160159951"),
    @("100091401", "AVR2109160159955"),
    @("100091403", "AVR2109160159962"),
    @("100091412", "AVR2109160159993"),
    @("100091413", "AVR2109160159995"),
    @("100091414", "AVR2109160159997")
)

$accountsData = @(
    @("PS_KimberlySteele890355@rep.com", "REP", "400780", "''No Cost Signup'"),
    @("PS_LarryCooper438458@rep.com", "REP", "400781", "''No Cost Signup'"),
    @("PS_CarlWillis987884@rep.com", "REP", "400782", "''Donation Signup'"),
    @("PS_DianeStark251159@rep.com", "REP", "400783", "''30`$ kit Signup'"),
    @("PS_MatthewDalton830636@cust.com", "Cust", "", "Attached"),
    @("PS_DouglasHale299767@cust.com", "Cust", "", "Unattached"),
    @("PS_JeremyJones818349@rep.com", "REP", "400788", "''Donation Signup'"),
    @("PS_ShaunThomas345399@rep.com", "REP", "400789", "''30`$ kit Signup'"),
    @("PS_MelindaSingh531223@rep.com", "REP", "400791", "''30`$ kit Signup'"),
    @("PS_TanyaLawson263831@cust.com", "Cust", "", "Unattached"),
    @("PS_TonyaMoore484276@cust.com", "Cust", "", "Unattached"),
    @("PS_ReneeRiley097838@cust.com", "Cust", "", "Unattached"),
    @("PS_StevenWagner596338@cust.com", "Cust", "", "Unattached"),
    @("PS_BrendaLevine471513@cust.com", "Cust", "", "Unattached"),
    @("PS_TiffanyWaters591843@cust.com", "Cust", "", "Unattached"),
    @("PS_JustinParks532136@cust.com", "Cust", "", "Unattached"),
    @("PS_MarkMerritt066755@cust.com", "Cust", "", "Unattached"),
    @("PS_LaurenDay052291@cust.com", "Cust", "", "Unattached"),
    @("PS_AndreHoffman032937@cust.com", "Cust", "", "Unattached"),
    @("PS_AlanHarris521913@rep.com", "REP", "400804", "''30`$ kit Signup'"),
    @("PS_AndrewStafford029758@rep.com", "REP", "400806", "''30`$ kit Signup'")
)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Order#" sheet: append rows 25-52 (Order No / Master Order No)
# ---------------------------------------------------------------------
$wsOrder = $wb.Worksheets.Item("Order#")
$startRow = 25
for ($i = 0; $i -lt $orderData.Count; $i++) {
    $r = $startRow + $i
    $pair = $orderData[$i]

    $cellA = $wsOrder.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $pair[0]

    $cellB = $wsOrder.Cells.Item($r, 2)
    $cellB.Value = $pair[1]
}

# ---------------------------------------------------------------------
# "Accounts" sheet: append rows 73-93 (Email / Type / Bee No. / Notes)
# ---------------------------------------------------------------------
$wsAccounts = $wb.Worksheets.Item("Accounts")
$startRow2 = 73
for ($i = 0; $i -lt $accountsData.Count; $i++) {
    $r = $startRow2 + $i
    $row = $accountsData[$i]

    $wsAccounts.Cells.Item($r, 1).Value = $row[0]
    $wsAccounts.Cells.Item($r, 2).Value = $row[1]

    if ($row[2] -ne "") {
        $cellC = $wsAccounts.Cells.Item($r, 3)
        $cellC.NumberFormat = "@"
        $cellC.Value = $row[2]
    }

    $wsAccounts.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------
# Column width tweaks on "Accounts" (matches target widths 62 / 6 / 32)
# ---------------------------------------------------------------------
$wsAccounts.Columns.Item(1).ColumnWidth = 61.16666666666667
$wsAccounts.Columns.Item(2).ColumnWidth = 5.166666666666667
$wsAccounts.Columns.Item(4).ColumnWidth = 31.166666666666668

# ---------------------------------------------------------------------
# Selection / scroll position on "Accounts" sheet view
# ---------------------------------------------------------------------
[void]$wsAccounts.Range("D77").Select()
